$d = $word.ActiveDocument

# --- 1) Title paragraph: "{{nombreInvestigacion}} ({{comite}})" / "Competent Body"
#     becomes "{{memoria.comite.nombre}} ({{memoria.comite.codigo}})" on its own
#     paragraph, followed by "Competent Body" on a new paragraph (the manual line
#     break between them is turned into a real paragraph break).

$d.Content.Find.Execute("nombreInvestigacion", $false, $false, $false, $false, $false, `
    $true, 1, $false, "memoria.comite.nombre", 2)

# The bookmark __DdeLink__74_3013164213 wraps just the word "comite"; replace its
# contents in place so the bookmark keeps wrapping the (now single-letter) "m".
$bmComite = $d.Bookmarks("__DdeLink__74_3013164213")
$bmComite.Range.Text = "m"

$d.Content.Find.Execute("}})", $false, $false, $false, $false, $false, `
    $true, 1, $false, "emoria.comite.codigo}})", 2)

# Turn the manual line break right after "}})" into a paragraph break, which
# splits "Competent Body" into its own paragraph.
$d.Content.Find.Execute("}})^l", $false, $false, $false, $false, $false, `
    $true, 1, $false, "}})^p", 2)

# --- 2) "{{nombreSecretario}}" -> "{{secretario.nombre}} {{secretario.apellidos}}"
#     split across two runs, separated by a new bookmark, and the pre-existing
#     (empty) bookmark moved to just before the first run.

$d.Content.Find.Execute("{{nombreSecretario}}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{{secretario.nombre}}", 2)

$rSecretario = $d.Content
$rSecretario.Find.Execute("{{secretario.nombre}}")
$secStart = $rSecretario.Start
$secEnd = $rSecretario.End

$insertApellidos = $d.Range($secEnd, $secEnd)
$insertApellidos.InsertAfter(" {{secretario.apellidos}}")

$newBmRange = $d.Range($secEnd, $secEnd)
$d.Bookmarks.Add("__DdeLink__38_20037553881", $newBmRange)

$oldBmRange = $d.Range($secStart, $secStart)
$d.Bookmarks.Add("__DdeLink__38_2003755388", $oldBmRange)

# --- 3) Simple placeholder renames

$d.Content.Find.Execute("tituloProyecto", $false, $false, $false, $false, $false, `
    $true, 1, $false, "peticionEvaluacion.titulo", 2)

$d.Content.Find.Execute(": {{nombreInvestigador}}", $false, $false, $false, $false, $false, `
    $true, 1, $false, ": {{investigador.nombre}} {{investigador.apellidos}}", 2)
